# Nightly cryptos-list refresh (prices & 1h volume %) + one name/URL
# ordering swap between FraxShare and Quant (rows 41-42).
#
# Price/percent cells are numeric-looking text (e.g. "217.83", "  +0.40%  ")
# stored as strings in the original sheet. Assigning such a string via
# Range.Value lets the COM layer auto-detect and coerce it to a number,
# so a leading apostrophe is used to force literal text entry -- exactly
# like typing '217.83 into Excel -- preserving the original text cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.128.01"
$ws.Range("E2").Value = "'  +0.29%  "
$ws.Range("D3").Value = "'1.654.52"
$ws.Range("E3").Value = "'  -0.14%  "
$ws.Range("E4").Value = "'  -0.16%  "
$ws.Range("D5").Value = "'217.83"
$ws.Range("E5").Value = "'  +0.40%  "
$ws.Range("D6").Value = "'0.5296"
$ws.Range("E6").Value = "'  +1.76%  "
$ws.Range("E7").Value = "'  -0.14%  "
$ws.Range("E8").Value = "'  -0.64%  "
$ws.Range("D9").Value = "'0.06327"
$ws.Range("E9").Value = "'  +1.10%  "
$ws.Range("E10").Value = "'  -0.62%  "
$ws.Range("D11").Value = "'0.07802"
$ws.Range("E11").Value = "'  +0.99%  "
$ws.Range("D12").Value = "'4.522"
$ws.Range("E12").Value = "'  +1.41%  "
$ws.Range("D13").Value = "'1.647.27"
$ws.Range("E13").Value = "'  -0.31%  "
$ws.Range("D14").Value = "'1.882.30"
$ws.Range("E14").Value = "'  -0.17%  "
$ws.Range("D15").Value = "'0.5485"
$ws.Range("E15").Value = "'  +1.06%  "
$ws.Range("D16").Value = "'0.0₅8216"
$ws.Range("E16").Value = "'  +1.90%  "
$ws.Range("D17").Value = "'65.42"
$ws.Range("E17").Value = "'  +1.15%  "
$ws.Range("D18").Value = "'26.124.98"
$ws.Range("E18").Value = "'  +0.17%  "
$ws.Range("E19").Value = "'  -0.02%  "
$ws.Range("D20").Value = "'4.598"
$ws.Range("E20").Value = "'  +0.56%  "
$ws.Range("D21").Value = "'190.99"
$ws.Range("E21").Value = "'  +0.16%  "
$ws.Range("E22").Value = "'  +0.89%  "
$ws.Range("D23").Value = "'6.027"
$ws.Range("E23").Value = "'  +1.13%  "
$ws.Range("E24").Value = "'  -0.23%  "
$ws.Range("D25").Value = "'145.22"
$ws.Range("E25").Value = "'  +5.39%  "
$ws.Range("D26").Value = "'0.1229"
$ws.Range("E26").Value = "'  -0.39%  "
$ws.Range("D27").Value = "'7.216"
$ws.Range("E27").Value = "'  -0.13%  "
$ws.Range("D28").Value = "'16.00"
$ws.Range("E28").Value = "'  -0.65%  "
$ws.Range("D29").Value = "'1.458"
$ws.Range("E29").Value = "'  +3.93%  "
$ws.Range("D30").Value = "'0.05784"
$ws.Range("E30").Value = "'  -1.96%  "
$ws.Range("D31").Value = "'1.274"
$ws.Range("E31").Value = "'  -0.17%  "
$ws.Range("D32").Value = "'3.552"
$ws.Range("E32").Value = "'  +0.98%  "
$ws.Range("D33").Value = "'3.270"
$ws.Range("E33").Value = "'  +0.66%  "
$ws.Range("D34").Value = "'1.600"
$ws.Range("E34").Value = "'  +2.45%  "
$ws.Range("D35").Value = "'2.801"
$ws.Range("E35").Value = "'  +1.38%  "
$ws.Range("D36").Value = "'0.9490"
$ws.Range("E36").Value = "'  -0.08%  "
$ws.Range("D37").Value = "'2.415"
$ws.Range("E37").Value = "'  -0.07%  "
$ws.Range("D38").Value = "'0.5757"
$ws.Range("E38").Value = "'  +2.11%  "
$ws.Range("D39").Value = "'0.01611"
$ws.Range("E39").Value = "'  +1.38%  "
$ws.Range("D40").Value = "'0.8564"
$ws.Range("E40").Value = "'  +1.30%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.772"
$ws.Range("E41").Value = "'  -1.94%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'104.56"
$ws.Range("E42").Value = "'  +3.81%  "
$ws.Range("E43").Value = "'  -0.04%  "
$ws.Range("D44").Value = "'1.034.14"
$ws.Range("E44").Value = "'  +3.45%  "
$ws.Range("D45").Value = "'1.796.05"
$ws.Range("E45").Value = "'  -0.29%  "
$ws.Range("E46").Value = "'  +0.83%  "
$ws.Range("D47").Value = "'1.009"
$ws.Range("E47").Value = "'  +0.86%  "
$ws.Range("D48").Value = "'0.4336"
$ws.Range("E48").Value = "'  +0.49%  "
$ws.Range("D49").Value = "'7.888"
$ws.Range("E49").Value = "'  -0.61%  "
$ws.Range("E50").Value = "'  -0.13%  "
$ws.Range("D51").Value = "'1.448"
$ws.Range("E51").Value = "'  -0.51%  "
